# Weekly symbol-list / price refresh (coinranking scrape), 2023-01-30 17:xx UTC
# Updates Price (D) and Volume(1h) (E) values for most rows, and re-sorts the
# CoinExToken..ZBToken block (rows 17-24) by one position, swapping in TigerCash
# at the top and cycling CoinExToken to the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''310.06'
$ws.Range("E2").Value = '''-2.55%'
$ws.Range("D3").Value = '''37.70'
$ws.Range("E3").Value = '''-4.87%'
$ws.Range("D4").Value = '''5.104'
$ws.Range("E4").Value = '''-0.69%'
$ws.Range("D5").Value = '''0.07856'
$ws.Range("D6").Value = '''1.964'
$ws.Range("E6").Value = '''-8.72%'
$ws.Range("D7").Value = '''4.370'
$ws.Range("E7").Value = '''1.93%'
$ws.Range("D8").Value = '''8.305'
$ws.Range("E8").Value = '''-0.01%'
$ws.Range("D9").Value = '''3.104'
$ws.Range("E9").Value = '''-6.29%'
$ws.Range("D10").Value = '''0.9275'
$ws.Range("E10").Value = '''-0.41%'
$ws.Range("D11").Value = '''0.1350'
$ws.Range("E11").Value = '''-4.41%'
$ws.Range("D12").Value = '''0.1967'
$ws.Range("E12").Value = '''-0.38%'
$ws.Range("D13").Value = '''0.08955'
$ws.Range("E13").Value = '''-1.16%'
$ws.Range("D14").Value = '''0.03476'
$ws.Range("E14").Value = '''-0.03%'
$ws.Range("D15").Value = '''0.09709'
$ws.Range("E15").Value = '''-0.94%'
$ws.Range("D16").Value = '''0.001394'
$ws.Range("E16").Value = '''-0.56%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.006010'
$ws.Range("E17").Value = '''-2.23%'
$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D18").Value = '''0.007506'
$ws.Range("E18").Value = '''1,778.35%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '''3.592'
$ws.Range("E19").Value = '''-2.51%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3465'
$ws.Range("E20").Value = '''-0.23%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1295'
$ws.Range("E21").Value = '''0.05%'
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").Value = '''5.004'
$ws.Range("E22").Value = '''2.13%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2513'
$ws.Range("E23").Value = '''2.63%'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = '''0.04348'
$ws.Range("E24").Value = '''0.48%'
$ws.Range("D25").Value = '''0.001225'
$ws.Range("E25").Value = '''-0.15%'
$ws.Range("D26").Value = '''0.004543'
$ws.Range("E26").Value = '''-4.73%'
$ws.Range("E27").Value = '''3.96%'
$ws.Range("D39").Value = '''0.02284'
$ws.Range("E39").Value = '''3.12%'
$ws.Range("D40").Value = '''0.05048'
$ws.Range("E40").Value = '''-3.46%'
$ws.Range("D41").Value = '''0.007610'
$ws.Range("E41").Value = '''1.12%'
$ws.Range("D42").Value = '''0.009784'
$ws.Range("E42").Value = '''1.06%'
$ws.Range("D43").Value = '''0.1357'
$ws.Range("E43").Value = '''-1.70%'
$ws.Range("D44").Value = '''0.002042'
$ws.Range("E44").Value = '''-3.41%'
$ws.Range("D45").Value = '''0.008785'
$ws.Range("E45").Value = '''-10.87%'
$ws.Range("D46").Value = '''0.00006810'
$ws.Range("E46").Value = '''3.37%'
$ws.Range("D47").Value = '''0.00000000751'
$ws.Range("E47").Value = '''0.12%'
$ws.Range("D48").Value = '''0.003003'
$ws.Range("E48").Value = '''8.61%'
$ws.Range("E49").Value = '''8.44%'
$ws.Range("D50").Value = '''0.00002102'
$ws.Range("E50").Value = '''0.12%'
$ws.Range("D51").Value = '''0.0002002'
$ws.Range("E51").Value = '''0.12%'
